$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1290.5385
$ws.Range("I40").Value = 760
$ws.Range("K40").Value = 760
$ws.Range("M40").Value = -585

$ws.Range("H113").Value = 90914300
$ws.Range("I113").Value = 166670300
$ws.Range("K113").Value = 166670300
$ws.Range("M113").Value = -166667046

$ws.Range("H129").Value = 182817.84
$ws.Range("I129").Value = 450
$ws.Range("J129").Value = 189699.64
$ws.Range("K129").Value = 1350
$ws.Range("L129").Value = 569098.92
$ws.Range("M129").Value = 3650
$ws.Range("N129").Value = -579098.92

$ws.Range("H138").Value = 2844.926
$ws.Range("I138").Value = 1946.6923
$ws.Range("J138").Value = 3679
$ws.Range("K138").Value = 5840.0769
$ws.Range("L138").Value = 11037
$ws.Range("M138").Value = -700.0769
$ws.Range("N138").Value = -21317

$ws.Range("H141").Value = 2634.625
$ws.Range("I141").Value = 2244.5
$ws.Range("K141").Value = 6733.5
$ws.Range("M141").Value = -1553.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 865.4211
$ws.Range("J2").Value = 585.3333
$ws.Range("L2").Value = 585.3333
$ws.Range("N2").Value = -811.3333

$ws.Range("H61").Value = 4921.875
$ws.Range("I61").Value = 5229.1665
$ws.Range("K61").Value = 5229.1665
$ws.Range("M61").Value = -5017.1665

$ws.Range("H74").Value = 28572964
$ws.Range("I74").Value = 43478896
$ws.Range("J74").Value = 3258.25
$ws.Range("K74").Value = 43478896
$ws.Range("L74").Value = 3258.25
$ws.Range("M74").Value = -43478022
$ws.Range("N74").Value = -5006.25

$ws.Range("H77").Value = 28572964
$ws.Range("I77").Value = 43478896
$ws.Range("J77").Value = 3258.25
$ws.Range("K77").Value = 217394480
$ws.Range("L77").Value = 16291.25
$ws.Range("M77").Value = -217390112
$ws.Range("N77").Value = -25027.25

$ws.Range("H110").Value = 830.125
$ws.Range("I110").Value = 748.7143
$ws.Range("J110").Value = 1400
$ws.Range("K110").Value = 748.7143
$ws.Range("L110").Value = 1400
$ws.Range("M110").Value = 1296.2857
$ws.Range("N110").Value = -5490

$ws.Range("H116").Value = 865.4211
$ws.Range("J116").Value = 585.3333
$ws.Range("L116").Value = 585.3333
$ws.Range("N116").Value = -5173.3333

$ws.Range("H122").Value = 2976.375
$ws.Range("I122").Value = 1973.1428
$ws.Range("K122").Value = 5919.428400000001
$ws.Range("M122").Value = -3469.428400000001

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 15742.737
$ws.Range("I132").Value = 2429
$ws.Range("K132").Value = 7287
$ws.Range("M132").Value = -4757

$ws.Range("H136").Value = 4921.875
$ws.Range("I136").Value = 5229.1665
$ws.Range("K136").Value = 15687.4995
$ws.Range("M136").Value = -13137.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 865.4211
$ws.Range("J3").Value = 585.3333
$ws.Range("L3").Value = 585.3333
$ws.Range("N3").Value = -813.3333

$ws.Range("H22").Value = 696.2857
$ws.Range("I22").Value = 644.8
$ws.Range("J22").Value = 825
$ws.Range("K22").Value = 644.8
$ws.Range("L22").Value = 825
$ws.Range("M22").Value = -471.8
$ws.Range("N22").Value = -1171

$ws.Range("H64").Value = 276.8125
$ws.Range("I64").Value = 129
$ws.Range("J64").Value = 424.625
$ws.Range("K64").Value = 129
$ws.Range("L64").Value = 424.625
$ws.Range("M64").Value = 96
$ws.Range("N64").Value = -874.625

$ws.Range("H67").Value = 276.8125
$ws.Range("I67").Value = 129
$ws.Range("J67").Value = 424.625
$ws.Range("K67").Value = 129
$ws.Range("L67").Value = 424.625
$ws.Range("M67").Value = 651
$ws.Range("N67").Value = -1984.625

$ws.Range("H105").Value = 1837.2128
$ws.Range("I105").Value = 1637.4166
$ws.Range("K105").Value = 1637.4166
$ws.Range("M105").Value = 109.5834

$ws.Range("H134").Value = 3465.1794
$ws.Range("I134").Value = 3560.75
$ws.Range("K134").Value = 10682.25
$ws.Range("M134").Value = -8147.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1401.375
$ws.Range("I16").Value = 1118.5
$ws.Range("J16").Value = 2250
$ws.Range("K16").Value = 1118.5
$ws.Range("L16").Value = 2250
$ws.Range("M16").Value = -831.5
$ws.Range("N16").Value = -2824

$ws.Range("H31").Value = 4141.7744
$ws.Range("I31").Value = 3614.3845
$ws.Range("K31").Value = 3614.3845
$ws.Range("M31").Value = -3319.3845

$ws.Range("H34").Value = 4141.7744
$ws.Range("I34").Value = 3614.3845
$ws.Range("K34").Value = 3614.3845
$ws.Range("M34").Value = -3412.3845

$ws.Range("H105").Value = 6251068.5
$ws.Range("I105").Value = 7813401
$ws.Range("J105").Value = 1737.75
$ws.Range("K105").Value = 7813401
$ws.Range("L105").Value = 1737.75
$ws.Range("M105").Value = -7811654
$ws.Range("N105").Value = -5231.75

$ws.Range("H113").Value = 1401.375
$ws.Range("I113").Value = 1118.5
$ws.Range("J113").Value = 2250
$ws.Range("K113").Value = 1118.5
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 1051.5
$ws.Range("N113").Value = -6590

$ws.Range("H132").Value = 4094.5
$ws.Range("I132").Value = 2954.5454
$ws.Range("K132").Value = 8863.636200000001
$ws.Range("M132").Value = -6333.636200000001

$ws.Range("H134").Value = 1417.2
$ws.Range("I134").Value = 1175.8
$ws.Range("K134").Value = 3527.4
$ws.Range("M134").Value = -992.3999999999996

$ws.Range("H141").Value = 13585.111
$ws.Range("J141").Value = 13585.111
$ws.Range("L141").Value = 13585.111
$ws.Range("N141").Value = -23945.111

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 378.14285
$ws.Range("I9").Value = 268.33334
$ws.Range("J9").Value = 460.5
$ws.Range("K9").Value = 805.0000200000001
$ws.Range("L9").Value = 1381.5
$ws.Range("M9").Value = -581.0000200000001
$ws.Range("N9").Value = -1829.5

$ws.Range("H131").Value = 698.2371000000001
$ws.Range("J131").Value = 698.2371000000001
$ws.Range("L131").Value = 2094.7113
$ws.Range("N131").Value = -12174.7113

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10792.533
$ws.Range("I70").Value = 4727.4287
$ws.Range("K70").Value = 4727.4287
$ws.Range("M70").Value = -4457.4287

$ws.Range("H73").Value = 10792.533
$ws.Range("I73").Value = 4727.4287
$ws.Range("K73").Value = 4727.4287
$ws.Range("M73").Value = -3791.4287

$ws.Range("H132").Value = 21145.932
$ws.Range("I132").Value = 4010.1738
$ws.Range("J132").Value = 86833
$ws.Range("K132").Value = 12030.5214
$ws.Range("L132").Value = 260499
$ws.Range("M132").Value = -9500.5214
$ws.Range("N132").Value = -265559

$ws.Range("H136").Value = 7966.1333
$ws.Range("J136").Value = 7966.1333
$ws.Range("L136").Value = 23898.3999
$ws.Range("N136").Value = -28998.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 144204.14
$ws.Range("I35").Value = 144204.14
$ws.Range("K35").Value = 144204.14
$ws.Range("M35").Value = -143868.14

$ws.Range("H55").Value = 841.6667
$ws.Range("I55").Value = 907.2727
$ws.Range("K55").Value = 907.2727
$ws.Range("M55").Value = -734.2727

$ws.Range("H61").Value = 2943.1
$ws.Range("I61").Value = 1315.4117
$ws.Range("K61").Value = 1315.4117
$ws.Range("M61").Value = -1113.4117

$ws.Range("H113").Value = 2943.1
$ws.Range("I113").Value = 1315.4117
$ws.Range("K113").Value = 1315.4117
$ws.Range("M113").Value = 854.5882999999999

$ws.Range("H122").Value = 1512118
$ws.Range("J122").Value = 4548.3335
$ws.Range("L122").Value = 13645.0005
$ws.Range("N122").Value = -18545.0005

$ws.Range("H132").Value = 483964.62
$ws.Range("I132").Value = 754396.2
$ws.Range("K132").Value = 2263188.6
$ws.Range("M132").Value = -2260658.6

$ws.Range("H136").Value = 1325.2759
$ws.Range("I136").Value = 1230.9584
$ws.Range("J136").Value = 1778
$ws.Range("K136").Value = 3692.8752
$ws.Range("L136").Value = 5334
$ws.Range("M136").Value = -1142.8752
$ws.Range("N136").Value = -10434

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1126842.2
$ws.Range("I113").Value = 821.8889
$ws.Range("J113").Value = 4504903.5
$ws.Range("K113").Value = 2465.6667
$ws.Range("L113").Value = 13514710.5
$ws.Range("M113").Value = -295.6667000000002
$ws.Range("N113").Value = -13519050.5

$ws.Range("H126").Value = 1519.125
$ws.Range("I126").Value = 1153.6086
$ws.Range("K126").Value = 3460.8258
$ws.Range("M126").Value = -990.8258000000001

$ws.Range("H132").Value = 1331.4375
$ws.Range("I132").Value = 950.2857
$ws.Range("K132").Value = 2850.8571
$ws.Range("M132").Value = -320.8571000000002

$ws.Range("H136").Value = 26471260
$ws.Range("I136").Value = 49156788
$ws.Range("J136").Value = 4810.5557
$ws.Range("K136").Value = 147470364
$ws.Range("L136").Value = 14431.6671
$ws.Range("M136").Value = -147467814
$ws.Range("N136").Value = -19531.6671
